$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 858.1111
$ws.Range("J129").Value = 961
$ws.Range("L129").Value = 2883
$ws.Range("N129").Value = -12883
$ws.Range("H137").Value = 3546.4043
$ws.Range("I137").Value = 2950.8948
$ws.Range("J137").Value = 6060.778
$ws.Range("K137").Value = 8852.6844
$ws.Range("L137").Value = 18182.334
$ws.Range("M137").Value = -6302.6844
$ws.Range("N137").Value = -23282.334
$ws.Range("H138").Value = 2411.9697
$ws.Range("I138").Value = 1212.2693
$ws.Range("J138").Value = 2839.2603
$ws.Range("K138").Value = 3636.8079
$ws.Range("L138").Value = 8517.7809
$ws.Range("M138").Value = 1503.1921
$ws.Range("N138").Value = -18797.7809
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10008.404
$ws.Range("I32").Value = 6474.491
$ws.Range("K32").Value = 6474.491
$ws.Range("M32").Value = -6187.491
$ws.Range("H76").Value = 23211.6
$ws.Range("J76").Value = 23211.6
$ws.Range("L76").Value = 23211.6
$ws.Range("N76").Value = -23887.6
$ws.Range("H79").Value = 23211.6
$ws.Range("J79").Value = 23211.6
$ws.Range("L79").Value = 23211.6
$ws.Range("N79").Value = -25551.6
$ws.Range("H107").Value = 30000
$ws.Range("J107").Value = 30000
$ws.Range("L107").Value = 30000
$ws.Range("N107").Value = -37680
$ws.Range("H109").Value = 30980
$ws.Range("J109").Value = 30980
$ws.Range("L109").Value = 30980
$ws.Range("N109").Value = -33754
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 350
$ws.Range("I64").Value = 200
$ws.Range("J64").Value = 500
$ws.Range("K64").Value = 200
$ws.Range("L64").Value = 500
$ws.Range("M64").Value = 25
$ws.Range("N64").Value = -950
$ws.Range("H67").Value = 350
$ws.Range("I67").Value = 200
$ws.Range("J67").Value = 500
$ws.Range("K67").Value = 200
$ws.Range("L67").Value = 500
$ws.Range("M67").Value = 580
$ws.Range("N67").Value = -2060
$ws.Range("H99").Value = 3697.3684
$ws.Range("I99").Value = 1165
$ws.Range("J99").Value = 4866.154
$ws.Range("K99").Value = 1165
$ws.Range("L99").Value = 4866.154
$ws.Range("M99").Value = 333
$ws.Range("N99").Value = -7862.154
$ws.Range("H118").Value = 29390
$ws.Range("J118").Value = 29390
$ws.Range("L118").Value = 29390
$ws.Range("N118").Value = -32704
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4585.6924
$ws.Range("I99").Value = 2714.2856
$ws.Range("J99").Value = 6769
$ws.Range("K99").Value = 2714.2856
$ws.Range("L99").Value = 6769
$ws.Range("M99").Value = -1216.2856
$ws.Range("N99").Value = -9765
$ws.Range("H122").Value = 2281.5386
$ws.Range("I122").Value = 1644.762
$ws.Range("J122").Value = 4956
$ws.Range("K122").Value = 4934.286
$ws.Range("L122").Value = 14868
$ws.Range("M122").Value = -2484.286
$ws.Range("N122").Value = -19768
$ws.Range("H123").Value = 37445.715
$ws.Range("J123").Value = 37445.715
$ws.Range("L123").Value = 37445.715
$ws.Range("N123").Value = -47245.715
$ws.Range("H126").Value = 4585.6924
$ws.Range("I126").Value = 2714.2856
$ws.Range("J126").Value = 6769
$ws.Range("K126").Value = 8142.8568
$ws.Range("L126").Value = 20307
$ws.Range("M126").Value = -5672.8568
$ws.Range("N126").Value = -25247
$ws.Range("H134").Value = 5511.8887
$ws.Range("I134").Value = 6112.55
$ws.Range("K134").Value = 18337.65
$ws.Range("M134").Value = -15802.65
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 518.4474
$ws.Range("I113").Value = 526.1818
$ws.Range("J113").Value = 507.8125
$ws.Range("K113").Value = 1578.5454
$ws.Range("L113").Value = 1523.4375
$ws.Range("M113").Value = 591.4546
$ws.Range("N113").Value = -5863.4375
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H97").Value = 2932.111
$ws.Range("I97").Value = 2398.1667
$ws.Range("J97").Value = 4000
$ws.Range("K97").Value = 2398.1667
$ws.Range("L97").Value = 4000
$ws.Range("M97").Value = -1902.1667
$ws.Range("N97").Value = -4992
$ws.Range("H112").Value = 28390
$ws.Range("J112").Value = 28390
$ws.Range("L112").Value = 28390
$ws.Range("N112").Value = -30606
$ws.Range("H113").Value = 1885.5
$ws.Range("I113").Value = 1460
$ws.Range("J113").Value = 4013
$ws.Range("K113").Value = 1460
$ws.Range("L113").Value = 4013
$ws.Range("M113").Value = 710
$ws.Range("N113").Value = -8353
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H126").Value = 4115.8623
$ws.Range("I126").Value = 2974.468
$ws.Range("J126").Value = 5457
$ws.Range("K126").Value = 8923.403999999999
$ws.Range("L126").Value = 16371
$ws.Range("M126").Value = -6453.403999999999
$ws.Range("N126").Value = -21311
$ws.Range("H132").Value = 2435.3928
$ws.Range("I132").Value = 1243.95
$ws.Range("J132").Value = 5414
$ws.Range("K132").Value = 3731.85
$ws.Range("L132").Value = 16242
$ws.Range("M132").Value = -1201.85
$ws.Range("N132").Value = -21302
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 43000
$ws.Range("J75").Value = 43000
$ws.Range("L75").Value = 43000
$ws.Range("N75").Value = -44872
$ws.Range("H78").Value = 43000
$ws.Range("J78").Value = 43000
$ws.Range("L78").Value = 129000
$ws.Range("N78").Value = -138360
$ws.Range("H100").Value = 2438.6
$ws.Range("I100").Value = 2058
$ws.Range("J100").Value = 3326.6667
$ws.Range("K100").Value = 2058
$ws.Range("L100").Value = 3326.6667
$ws.Range("M100").Value = -1517
$ws.Range("N100").Value = -4408.6667
$ws.Range("H122").Value = 3788.4
$ws.Range("I122").Value = 3022.3872
$ws.Range("K122").Value = 9067.161599999999
$ws.Range("M122").Value = -6617.161599999999
$ws.Range("H133").Value = 54743.332
$ws.Range("J133").Value = 54743.332
$ws.Range("L133").Value = 54743.332
$ws.Range("N133").Value = -59803.332
